$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new functions were added with manual calculation mode turned on.
$excel.Calculation = -4135

# Row 22: verifyWarningIsDisplayedForTheField
$ws.Range("A22").Value = 21
$ws.Range("C22").Value = "verifyWarningIsDisplayedForTheField"
$ws.Range("B22").Value = "Text Field, DropDowns"
$ws.Range("D22").Value = "Accepts one parameter @locator. Checks if the warning symbol is displayed for the input fields"
$ws.Range("D22").WrapText = $true

# Row 23: verifyNoWarningIsDisplayedForTheField
$ws.Range("A23").Value = 22
$ws.Range("C23").Value = "verifyNoWarningIsDisplayedForTheField"
$ws.Range("B23").Value = "Text Field, DropDowns"
$ws.Range("D23").Value = "Accepts one parameter @locator. Verifies no warning symbol is displayed for the input fields"
$ws.Range("D23").WrapText = $true

# The Field Type column for the Button/Text Field rows picked up an
# explicit (re-applied) font while editing.
$ws.Range("B20:B23").Font.ThemeColor = 1

# Leave the selection where the editor ended up.
$ws.Range("D24").Select()
